$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of A2 and A3, while also trimming the "blank" string
# from four spaces to three spaces (matching the target diff).
$ws.Range("A2").Value = "   "
$ws.Range("A3").Value = "plainaddress"

# Update the selection shown in the sheet view to D6.
$ws.Range("D6").Select()
